$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(13.38702964782715, -17.40622901916504, -7.910804748535156),
    @(5.027063369750977, -9.992935180664062, -3.668930053710938),
    @(3.064489126205444, -10.04933261871338, -9.655683517456056),
    @(-11.23465538024902, -9.59632682800293, -12.52409744262695),
    @(0.3307132720947265, -10.98133277893066, -7.000489711761475),
    @(15.33979511260986, -22.00781440734864, -10.50343608856201),
    @(2.999754428863525, -15.48295783996582, 6.745935916900635),
    @(-6.992471218109131, 1.269168853759766, -5.262283325195312),
    @(1.578595161437989, -59.7728385925293, 11.89765167236328),
    @(-58.78142929077149, -6.325687885284424, -18.43083572387696),
    @(1.685664176940918, -10.93907737731934, -4.804312705993652),
    @(20.54140090942383, -30.76116371154785, -21.73070907592773),
    @(-8.643145561218262, -12.25782108306885, 13.20825862884522),
    @(-1.923628926277161, -1.592215418815613, -5.497756004333496),
    @(23.09793281555176, -65.25662231445312, 1.447998523712158),
    @(13.40246772766113, -2.679043769836426, -17.18461418151855),
    @(-13.39319038391113, -15.08681869506836, -24.70808029174805),
    @(47.38216018676758, -81.48814392089844, 17.70297622680664),
    @(-2.278211832046509, 0.415550947189331, -10.76592063903809),
    @(13.07781982421875, -39.51717376708984, 42.08852005004883),
    @(-72.39402770996094, -27.50520896911621, -27.84894561767578),
    @(-8.723164558410645, -8.500687599182129, -15.66308975219727),
    @(-7.880284786224365, -21.83874320983887, -8.497885704040527),
    @(-9.739827156066896, -0.5102891325950623, -2.436953544616699),
    @(-2.781617164611816, 10.22575092315674, -2.437583446502685),
    @(34.05283355712891, -58.47340393066406, 20.97058296203613),
    @(-11.52477169036865, 17.23179626464844, -1.970066547393799),
    @(-24.07291412353516, -20.21802139282227, -17.61478233337402),
    @(18.03611946105957, -66.51955413818359, 48.28945922851562),
    @(12.77012634277344, 6.439512729644775, -10.52412605285644)
)

$rowCount = $data.Length
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
